# Append the 2021 row (row 13) of copyright-export data to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous year row (A12, style index "1":
# bold, bordered, centered) onto the new row's year cell so the new
# row matches the look of every other year label.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's values.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 1025
$ws.Range("C13").Value = 6238
$ws.Range("D13").Value = 233
$ws.Range("E13").Value = 1005
$ws.Range("F13").Value = 435
$ws.Range("G13").Value = 668
$ws.Range("H13").Value = 371
$ws.Range("I13").Value = 166
$ws.Range("J13").Value = 58
$ws.Range("K13").Value = 887
$ws.Range("L13").Value = 446
$ws.Range("M13").Value = 512
$ws.Range("N13").Value = 726
$ws.Range("O13").Value = 12770
